# Commit: "Fruta / hortaliza, semanal"
# The underlying data table (Apio / Vega Monumental Concepción) gains one
# new weekly observation. A new row is inserted at sheet row 163, pushing
# every existing record (old rows 163-288) down by one (new rows 164-289).
# The dimension grows from A1:R288 to A1:R289.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 163; everything below shifts down one.
$ws.Rows("163:163").Insert()

# Populate the newly inserted row 163 with the new weekly record.
$ws.Range("A163").Value = 11
$ws.Range("B163").Value = 'Vega Monumental Concepción'
$ws.Range("C163").Value = 'Bíobío'
$ws.Range("D163").Value = 44741
$ws.Range("E163").Value = 8
$ws.Range("F163").Value = 100112017
$ws.Range("G163").Value = 'Apio'
$ws.Range("H163").Value = 'Americana (o)'
$ws.Range("I163").Value = 'Primera'
$ws.Range("J163").Value = 100
$ws.Range("K163").Value = 6000
$ws.Range("L163").Value = 7000
$ws.Range("M163").Value = 6500
$ws.Range("N163").Value = '$/docena de matas'
$ws.Range("O163").Value = 'Región de Coquimbo'
$ws.Range("P163").Value = 1083
$ws.Range("Q163").Value = 6
$ws.Range("R163").Value = 'Hortaliza'
